$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row at position 56 (pushes old rows 56-88 down to 57-89,
# extending the used range from A1:R88 to A1:R89).
$ws.Rows.Item(56).Insert()

# Populate the newly inserted row with the new weekly price record.
$ws.Cells.Item(56,1).Value = 8
$ws.Cells.Item(56,2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(56,3).Value = "Coquimbo"
$ws.Cells.Item(56,4).Value = 44488
$ws.Cells.Item(56,5).Value = 4
$ws.Cells.Item(56,6).Value = 100112044
$ws.Cells.Item(56,7).Value = "Perejil"
$ws.Cells.Item(56,8).Value = "Sin especificar"
$ws.Cells.Item(56,9).Value = "Primera"
$ws.Cells.Item(56,10).Value = 3000
$ws.Cells.Item(56,11).Value = 1300
$ws.Cells.Item(56,12).Value = 1500
$ws.Cells.Item(56,13).Value = 1400
$ws.Cells.Item(56,14).Value = "`$/atado 1 a 1,5 kilos"
$ws.Cells.Item(56,15).Value = "Provincia del Elquí"
$ws.Cells.Item(56,16).Value = 933
$ws.Cells.Item(56,17).Value = 1.5
$ws.Cells.Item(56,18).Value = "Hortaliza"
